$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refresh the indicator values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 213666.68000000005
$metrics.Range("B3").Value  = 187746.13
$metrics.Range("B4").Value  = 66144.350000000006
$metrics.Range("B5").Value  = 8743
$metrics.Range("B6").Value  = 5009912.4300000016
$metrics.Range("B7").Value  = 4229822.8100000005
$metrics.Range("B8").Value  = 1473104.18
$metrics.Range("B9").Value  = 194950
$metrics.Range("B10").Value = 33475293.420000009
$metrics.Range("B11").Value = 31505097.969999999
$metrics.Range("B12").Value = 11754826.220000001
$metrics.Range("B13").Value = 1292580

# Move the selection on the Metrics sheet to D12.
$metrics.Activate() | Out-Null
$metrics.Range("D12").Select() | Out-Null

# --- today sheet: formulas pull from Metrics automatically on recalc; ---
# --- just reposition the selection, keeping this sheet the active tab. ---
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("H14").Select() | Out-Null
